# Update database: shift the yearly columns forward by one fiscal year
# (1397..1401 instead of 1396..1400), update the publish-date row, and
# reset all financial figures to 0 pending the new read_price algorithm
# (per commit message "update database and change read_price algorithm").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")

# --- Row 8: fiscal-period column headers -------------------------------
$ws.Range("D8").Value = "12 ماهه منتهی به 1397/12"
$ws.Range("E8").Value = "12 ماهه منتهی به 1398/12"
$ws.Range("F8").Value = "12 ماهه منتهی به 1399/12"
$ws.Range("G8").Value = "12 ماهه منتهی به 1400/12"
$ws.Range("H8").Value = "12 ماهه منتهی به 1401/12"

# --- Row 9: publish dates ------------------------------------------------
$ws.Range("D9").Value = "1399-02-15 (10)"
$ws.Range("E9").Value = "1400-02-19 (11)"
$ws.Range("F9").Value = "1401-02-25 (12)"
$ws.Range("G9").Value = "1402-02-27 (12)"
$ws.Range("H9").Value = "1402-02-27 (3)"

# --- Row 11: فروش (Sales) -------------------------------------------------
$ws.Range("D11:H11").Value = 0

# --- Row 12 ---------------------------------------------------------------
$ws.Range("D12:H12").Value = 0

# --- Row 13 ---------------------------------------------------------------
$ws.Range("D13:H13").Value = 0

# --- Row 14 ---------------------------------------------------------------
$ws.Range("D14:H14").Value = 0

# --- Row 15 (D15 already "-") ---------------------------------------------
$ws.Range("E15").Value = "-"
$ws.Range("F15").Value = "-"
$ws.Range("G15").Value = "-"
$ws.Range("H15").Value = "-"

# --- Row 16 -----------------------------------------------------------------
$ws.Range("D16:H16").Value = 0

# --- Row 17 -----------------------------------------------------------------
$ws.Range("D17:H17").Value = 0

# --- Row 18 -----------------------------------------------------------------
$ws.Range("D18:H18").Value = 0

# --- Row 19 -----------------------------------------------------------------
$ws.Range("D19:H19").Value = 0

# --- Row 20 -----------------------------------------------------------------
$ws.Range("D20:H20").Value = 0

# --- Row 21 -----------------------------------------------------------------
$ws.Range("D21:H21").Value = 0

# --- Row 22 -----------------------------------------------------------------
$ws.Range("D22:H22").Value = 0

# --- Row 23 (all become "-") -------------------------------------------------
$ws.Range("D23").Value = "-"
$ws.Range("E23").Value = "-"
$ws.Range("F23").Value = "-"
$ws.Range("G23").Value = "-"
$ws.Range("H23").Value = "-"

# --- Row 24 -----------------------------------------------------------------
$ws.Range("D24:H24").Value = 0

# --- Row 25 (E25 becomes "-", F25 stays 0) -----------------------------------
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = "-"
$ws.Range("F25").Value = 0
$ws.Range("G25").Value = 0
$ws.Range("H25").Value = 0

# --- Row 26 -----------------------------------------------------------------
$ws.Range("D26:H26").Value = 0

# --- Row 27 -----------------------------------------------------------------
$ws.Range("D27:H27").Value = 0
